$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.907.71'
$ws.Range("E2").Value = '  -0.79%  '

$ws.Range("D3").Value = '2.213.81'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '241.74'
$ws.Range("E5").Value = '  -1.96%  '

$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  -0.55%  '

$ws.Range("D7").Value = '72.84'
$ws.Range("E7").Value = '  -2.31%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").Value = '0.605'
$ws.Range("E9").Value = '  -2.06%  '

$ws.Range("D10").Value = '42.27'
$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("D11").Value = '0.0955'
$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("D12").Value = '7.04'
$ws.Range("E12").Value = '  -1.71%  '

$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  +0.63%  '

$ws.Range("D14").Value = '2.547.76'
$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("D15").Value = '14.26'
$ws.Range("E15").Value = '  -1.52%  '

$ws.Range("D16").Value = '0.836'
$ws.Range("E16").Value = '  -1.98%  '

$ws.Range("D17").Value = '2.210.63'
$ws.Range("E17").Value = '  -1.60%  '

$ws.Range("D18").Value = '41.840.10'

$ws.Range("D19").Value = '0.0000107'
$ws.Range("E19").Value = '  +6.40%  '

$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").Value = '72.71'
$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.18'
$ws.Range("E21").Value = '  +0.60%  '

$ws.Range("D22").Value = '10.83'
$ws.Range("E22").Value = '  +20.82%  '

$ws.Range("D23").Value = '230.29'
$ws.Range("E23").Value = '  -0.60%  '

$ws.Range("D24").Value = '2.07'
$ws.Range("E24").Value = '  -6.52%  '

$ws.Range("E25").Value = '  +3.08%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("E27").Value = '  +1.04%  '

$ws.Range("E28").Value = '  -1.69%  '

$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("D30").Value = '168.17'
$ws.Range("E30").Value = '  -0.34%  '

$ws.Range("E31").Value = '  -1.10%  '

$ws.Range("E32").Value = '  +7.68%  '

$ws.Range("D33").Value = '0.0794'
$ws.Range("E33").Value = '  -3.40%  '

$ws.Range("D34").Value = '29.71'
$ws.Range("E34").Value = '  -3.04%  '

$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("E36").Value = '  -10.72%  '

$ws.Range("E37").Value = '  -4.11%  '

$ws.Range("E38").Value = '  -4.43%  '

$ws.Range("D39").Value = '14.00'
$ws.Range("E39").Value = '  +1.59%  '

$ws.Range("D40").Value = '65.58'
$ws.Range("E40").Value = '  +4.71%  '

$ws.Range("E41").Value = '  -2.44%  '

$ws.Range("E42").Value = '  -2.70%  '

$ws.Range("D43").Value = '0.197'
$ws.Range("E43").Value = '  -3.10%  '

$ws.Range("E44").Value = '  +0.88%  '

$ws.Range("E45").Value = '  -2.43%  '

$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("D47").Value = '2.44'
$ws.Range("E47").Value = '  +5.99%  '

$ws.Range("E48").Value = '  -0.53%  '

$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("D51").Value = '2.420.74'
$ws.Range("E51").Value = '  -1.68%  '
